# Shift provider numbers in column A from provider11..provider21
# to provider21..provider31 (rows 1-11 of the active sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 0; $i -lt 11; $i++) {
    $row = $i + 1
    $newNum = 21 + $i
    $ws.Cells.Item($row, 1).Value = "provider$newNum"
}
